$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 306, pushing the existing rows 306-330 down to 307-331.
$ws.Rows.Item(306).Insert()

# Populate the new row 306 with the new weekly price record.
$ws.Range("A306").Value = 9
$ws.Range("B306").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C306").Value = "Metropolitana"
$ws.Range("D306").Value = 44578
$ws.Range("E306").Value = 13
$ws.Range("F306").Value = 100112039
$ws.Range("G306").Value = "Ciboulette"
$ws.Range("H306").Value = "Sin especificar"
$ws.Range("I306").Value = "Primera"
$ws.Range("J306").Value = 160
$ws.Range("K306").Value = 1000
$ws.Range("L306").Value = 1200
$ws.Range("M306").Value = 1100
$ws.Range("N306").Value = "$/docena de atados"
$ws.Range("O306").Value = "Región Metropolitana"
$ws.Range("P306").Value = 367
$ws.Range("Q306").Value = 3
$ws.Range("R306").Value = "Hortaliza"
